$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 361

# Mapping of the shorthand correct_ans codes to full words.
$ansMap = @{ "y" = "left"; "b" = "center"; "r" = "right" }

for ($r = 2; $r -le $lastRow; $r++) {
    # Columns A-D hold image filenames like "face//face_NN.jpg"; rename the
    # "face" category to "book" wherever it occurs, keeping the rest intact.
    foreach ($col in @("A", "B", "C", "D")) {
        $cell = $ws.Range("$col$r")
        $val = $cell.Value()
        if ($val -and $val.ToString().StartsWith("face//")) {
            $cell.Value = "book//" + $val.ToString().Substring(6)
        }
    }

    # Column L holds single-letter answer codes; expand them to full words.
    $lCell = $ws.Range("L$r")
    $lVal = $lCell.Value()
    if ($lVal -and $ansMap.ContainsKey($lVal.ToString())) {
        $lCell.Value = $ansMap[$lVal.ToString()]
    }
}
